# Update concise_ms csv pattern
# Row 11 ("Marking"): B11 4 -> 5, C11 -1 -> -1.2
# Row 12 ("Total"):   B12 92 -> 115, C12 -1 -> -1.2, E12 "91/112" -> "113.8/140"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

$ws.Range("B12").Value = 115
$ws.Range("C12").Value = -1.2
$ws.Range("E12").Value = "113.8/140"
